# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (col I) and DialogAct (col J)
# values for the rows affected by the re-annotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;   I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 6;   I = "aa"; J = "Agree/Accept" }
    @{ Row = 26;  I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 32;  I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 37;  I = "%";  J = "Uninterpretable" }
    @{ Row = 39;  I = "aa"; J = "Agree/Accept" }
    @{ Row = 40;  I = "aa"; J = "Agree/Accept" }
    @{ Row = 46;  I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 61;  I = "aa"; J = "Agree/Accept" }
    @{ Row = 63;  I = "aa"; J = "Agree/Accept" }
    @{ Row = 71;  I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 86;  I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 107; I = "aa"; J = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
